$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data (row 23), following the same column layout as row 22.
$ws.Range("A23").Value = 131289488
$ws.Range("B23").Value = 57881
$ws.Range("D23").Value = "NT"
$ws.Range("E23").Value = 100049
$ws.Range("F23").Value = "Spillkråka"
$ws.Range("G23").Value = "Dryocopus martius"
$ws.Range("H23").Value = "(Linnaeus, 1758)"
# Blank-but-present text cells (mirrors the source export's empty string cells).
$ws.Range("I23").Value = "'"
$ws.Range("K23").Value = "'"
$ws.Range("L23").Value = "'"
$ws.Range("M23").Value = "gammalt bo"
$ws.Range("N23").Value = "'"
$ws.Range("P23").Value = "Skogen norr om Sjöberga, Ög"
$ws.Range("Q23").Value = 567462
$ws.Range("R23").Value = 6509761
$ws.Range("S23").Value = 10
$ws.Range("T23").Value = "Östergötland"
$ws.Range("U23").Value = "Norrköping"
$ws.Range("V23").Value = "Östergötland"
$ws.Range("W23").Value = "Simonstorp"
# Dates are stored as literal text in this sheet, not Excel date serials.
$ws.Range("Y23").Value = "'2026-02-21"
$ws.Range("AA23").Value = "'2026-02-21"
$ws.Range("AD23").Value = $false
$ws.Range("AE23").Value = $false
$ws.Range("AG23").Value = $false
$ws.Range("AT23").Value = "'"
$ws.Range("AW23").Value = "Anette Källman"
$ws.Range("AX23").Value = "Anette Källman"
$ws.Range("AY23").Value = "'"

Write-Output "Row 23 written"
